$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.452.11'
$ws.Cells.Item(2, 5).Value = '  +0.80%  '
$ws.Cells.Item(3, 4).Value = '2.301.07'
$ws.Cells.Item(3, 5).Value = '  +0.27%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '316.50'
$ws.Cells.Item(5, 5).Value = '  +1.25%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '102.63'
$ws.Cells.Item(6, 5).Value = '  -1.89%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.630'
$ws.Cells.Item(7, 5).Value = '  +0.64%  '
$ws.Cells.Item(8, 5).Value = '  +0.11%  '
$ws.Cells.Item(9, 5).Value = '  -0.19%  '
$ws.Cells.Item(10, 5).Value = '  -2.00%  '
$ws.Cells.Item(11, 5).Value = '  -0.29%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '8.39'
$ws.Cells.Item(12, 5).Value = '  +1.54%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.107'
$ws.Cells.Item(13, 5).Value = '  +0.47%  '
$ws.Cells.Item(14, 5).Value = '  -0.94%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '15.21'
$ws.Cells.Item(15, 5).Value = '  -1.92%  '
$ws.Cells.Item(16, 4).Value = '2.648.12'
$ws.Cells.Item(16, 5).Value = '  +0.19%  '
$ws.Cells.Item(17, 4).Value = '2.305.67'
$ws.Cells.Item(17, 5).Value = '  +0.43%  '
$ws.Cells.Item(18, 4).Value = '42.404.28'
$ws.Cells.Item(18, 5).Value = '  +0.84%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '7.43'
$ws.Cells.Item(19, 5).Value = '  -2.11%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0000105'
$ws.Cells.Item(20, 5).Value = '  +0.74%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '73.40'
$ws.Cells.Item(21, 5).Value = '  -1.25%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '3.55'
$ws.Cells.Item(22, 5).Value = '  +3.03%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '276.95'
$ws.Cells.Item(23, 5).Value = '  +7.80%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '11.36'
$ws.Cells.Item(24, 5).Value = '  +22.24%  '
$ws.Cells.Item(25, 5).Value = '  -1.63%  '
$ws.Cells.Item(26, 5).Value = '  -0.37%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.82'
$ws.Cells.Item(27, 5).Value = '  -1.21%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.33'
$ws.Cells.Item(28, 5).Value = '  +2.20%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '22.71'
$ws.Cells.Item(29, 5).Value = '  -0.06%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '37.18'
$ws.Cells.Item(30, 5).Value = '  +4.06%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '165.76'
$ws.Cells.Item(31, 5).Value = '  -0.11%  '
$ws.Cells.Item(32, 5).Value = '  -2.26%  '
$ws.Cells.Item(33, 5).Value = '  +1.11%  '
$ws.Cells.Item(34, 5).Value = '  +4.75%  '
$ws.Cells.Item(35, 2).Value = 'WEMIXToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.61'
$ws.Cells.Item(35, 5).Value = '  -10.35%  '
$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.118'
$ws.Cells.Item(36, 5).Value = '  -0.26%  '
$ws.Cells.Item(37, 5).Value = '  +3.86%  '
$ws.Cells.Item(38, 5).Value = '  +0.54%  '
$ws.Cells.Item(39, 5).Value = '  +2.03%  '
$ws.Cells.Item(40, 5).Value = '  +0.18%  '
$ws.Cells.Item(41, 5).Value = '  +3.45%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '69.60'
$ws.Cells.Item(42, 5).Value = '  -3.14%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '94.66'
$ws.Cells.Item(43, 5).Value = '  -3.73%  '
$ws.Cells.Item(44, 5).Value = '  -0.39%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '80.93'
$ws.Cells.Item(46, 5).Value = '  +8.85%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '12.00'
$ws.Cells.Item(47, 5).Value = '  -2.16%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '112.62'
$ws.Cells.Item(48, 5).Value = '  +0.56%  '
$ws.Cells.Item(49, 5).Value = '  -1.00%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '5.24'
$ws.Cells.Item(50, 5).Value = '  -1.56%  '
$ws.Cells.Item(51, 4).Value = '1.589.10'
$ws.Cells.Item(51, 5).Value = '  +1.76%  '
